$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format temporarily so numeric-looking price strings
# (e.g. "1.00", "0.999") are preserved exactly instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.873.98"
$ws.Range("E2").Value = "  +0.64%  "

$ws.Range("D3").Value = "3.328.98"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("D5").Value = "401.16"
$ws.Range("E5").Value = "  -1.98%  "

$ws.Range("D6").Value = "126.23"
$ws.Range("E6").Value = "  +10.38%  "

$ws.Range("D7").Value = "0.592"
$ws.Range("E7").Value = "  +3.67%  "

$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "0.658"
$ws.Range("E9").Value = "  +5.64%  "

$ws.Range("E10").Value = "  +3.29%  "

$ws.Range("D11").Value = "41.17"
$ws.Range("E11").Value = "  +3.87%  "

$ws.Range("E12").Value = "  -0.89%  "

$ws.Range("D13").Value = "3.880.76"
$ws.Range("E13").Value = "  +1.11%  "

$ws.Range("D14").Value = "8.31"
$ws.Range("E14").Value = "  +1.87%  "

$ws.Range("D15").Value = "19.32"
$ws.Range("E15").Value = "  +1.58%  "

$ws.Range("D16").Value = "3.338.99"
$ws.Range("E16").Value = "  +0.97%  "

$ws.Range("D17").Value = "60.970.36"
$ws.Range("E17").Value = "  +1.01%  "

$ws.Range("D18").Value = "11.27"
$ws.Range("E18").Value = "  +4.65%  "

$ws.Range("E19").Value = "  +0.95%  "

$ws.Range("D20").Value = "0.0000128"
$ws.Range("E20").Value = "  +11.76%  "

$ws.Range("D21").Value = "3.20"
$ws.Range("E21").Value = "  -4.22%  "

$ws.Range("D22").Value = "80.55"
$ws.Range("E22").Value = "  +9.09%  "

$ws.Range("D23").Value = "12.88"
$ws.Range("E23").Value = "  +4.07%  "

$ws.Range("D24").Value = "300.81"
$ws.Range("E24").Value = "  +1.87%  "

$ws.Range("D25").Value = "3.11"
$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("E26").Value = "  +10.31%  "

$ws.Range("D27").Value = "8.36"
$ws.Range("E27").Value = "  +12.00%  "

$ws.Range("D28").Value = "29.01"
$ws.Range("E28").Value = "  -0.45%  "

$ws.Range("D29").Value = "7.40"
$ws.Range("E29").Value = "  -2.34%  "

$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").Value = "0.113"
$ws.Range("E31").Value = "  +0.06%  "

$ws.Range("D32").Value = "11.43"
$ws.Range("E32").Value = "  +1.94%  "

$ws.Range("D33").Value = "2.54"
$ws.Range("E33").Value = "  +2.56%  "

$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("D35").Value = "40.86"
$ws.Range("E35").Value = "  +0.40%  "

$ws.Range("D36").Value = "0.0477"
$ws.Range("E36").Value = "  -2.60%  "

$ws.Range("D37").Value = "52.10"
$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("E38").Value = "  +0.36%  "

$ws.Range("D39").Value = "3.36"
$ws.Range("E39").Value = "  +0.43%  "

$ws.Range("D40").Value = "2.90"
$ws.Range("E40").Value = "  -4.68%  "

$ws.Range("D41").Value = "1.97"
$ws.Range("E41").Value = "  +4.66%  "

$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "135.22"
$ws.Range("E42").Value = "  +0.93%  "

$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "0.123"
$ws.Range("E43").Value = "  +2.68%  "

$ws.Range("D44").Value = "3.88"
$ws.Range("E44").Value = "  +1.85%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.279"
$ws.Range("E45").Value = "  -3.89%  "

$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Value = "16.66"
$ws.Range("E46").Value = "  +2.84%  "

$ws.Range("E47").Value = "  +1.41%  "

$ws.Range("D48").Value = "21.38"
$ws.Range("E48").Value = "  +2.08%  "

$ws.Range("D49").Value = "2.119.02"
$ws.Range("E49").Value = "  -0.67%  "

$ws.Range("D50").Value = "3.668.95"
$ws.Range("E50").Value = "  +0.97%  "

$ws.Range("E51").Value = "  -1.27%  "

# Restore column D back to the default (Normal) style/format.
$ws.Range("D2:D51").Style = "Normal"